$wb = $excel.ActiveWorkbook

# Worksheets (by name, to be robust to ordering)
$wsPackages = $wb.Worksheets.Item("packages")
$wsEntities = $wb.Worksheets.Item("entities")
$wsAttributes = $wb.Worksheets.Item("attributes")
$wsErrors = $wb.Worksheets.Item("solverdportal_experiment_errors")

# 1) Bump the package description to the new version string (v0.9 -> v0.91)
$wsPackages.Range("C2").Value = "New official datatransfers from GPAP, CNAG (v0.91)"

# 2) Register the new "experiment_counts" entity in the entities table
$wsEntities.Range("A4").Value = "experiment_counts"
$wsEntities.Range("B4").Value = "GPAP Experiment Counts"
$wsEntities.Range("C4").Value = "Summary of errors"
$wsEntities.Range("D4").Value = "solverdportal"

# 3) Add the attribute definitions for the new solverdportal_experiment_counts entity
$wsAttributes.Range("A23").Value = "solverdportal_experiment_counts"
$wsAttributes.Range("B23").Value = "has_error"
$wsAttributes.Range("D23").Value = "bool"

$wsAttributes.Range("A24").Value = "solverdportal_experiment_counts"
$wsAttributes.Range("B24").Value = "error_type"
$wsAttributes.Range("D24").Value = "string"
$wsAttributes.Range("E24").Value = $true

$wsAttributes.Range("A25").Value = "solverdportal_experiment_counts"
$wsAttributes.Range("B25").Value = "count"
$wsAttributes.Range("D25").Value = "int"

# 4) Update sheet views / selections / active tab to match the saved workbook state

# packages becomes the active tab with C3 selected
$wsPackages.Activate() | Out-Null
$wsPackages.Range("C3").Select() | Out-Null

# entities: selection moves to A5
$wsEntities.Activate() | Out-Null
$wsEntities.Range("A5").Select() | Out-Null

# attributes: scrolled so row 2 is at the top, with A25 selected
$wsAttributes.Activate() | Out-Null
$win = $wb.Windows.Item(1)
$wsAttributes.Range("A25").Select() | Out-Null
$win.ScrollRow = 2
$win.ScrollColumn = 1

# solverdportal_experiment_errors: selection stays at A11 (unchanged), no longer the active tab
$wsErrors.Range("A11").Select() | Out-Null

# Re-activate packages last so it ends up as the active sheet/tab
$wsPackages.Activate() | Out-Null
